# chore: adapt column header formatting to respective input file names
#
# The sheet holds a diff table comparing two AHB format versions (FV2410 "old"
# vs FV2504 "new"). Headers used generic "_old"/"_new" suffixes; rename them
# to the concrete format-version suffixes, wrap the data range in a proper
# Excel Table, and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename the header cells: "_old" -> "_FV2410", "_new" -> "_FV2504" ---
$ws.Range("A1").Value = "Segmentname_FV2410"
$ws.Range("B1").Value = "Segmentgruppe_FV2410"
$ws.Range("C1").Value = "Segment_FV2410"
$ws.Range("D1").Value = "Datenelement_FV2410"
$ws.Range("E1").Value = "Segment ID_FV2410"
$ws.Range("F1").Value = "Code_FV2410"
$ws.Range("G1").Value = "Qualifier_FV2410"
$ws.Range("H1").Value = "Beschreibung_FV2410"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2410"
$ws.Range("J1").Value = "Bedingung_FV2410"
$ws.Range("K1").Value = "diff"
$ws.Range("L1").Value = "Segmentname_FV2504"
$ws.Range("M1").Value = "Segmentgruppe_FV2504"
$ws.Range("N1").Value = "Segment_FV2504"
$ws.Range("O1").Value = "Datenelement_FV2504"
$ws.Range("P1").Value = "Segment ID_FV2504"
$ws.Range("Q1").Value = "Code_FV2504"
$ws.Range("R1").Value = "Qualifier_FV2504"
$ws.Range("S1").Value = "Beschreibung_FV2504"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2504"
$ws.Range("U1").Value = "Bedingung_FV2504"

# --- 2) Wrap A1:U79 in an Excel Table, keeping the original header styling
#        (bold/filled/bordered cells) instead of the dxf Excel auto-creates
#        to "preserve" pre-existing header formatting, and without applying
#        a named table style. ---
$headerRange = $ws.Range("A1:U1")
$styleBackup = $ws.Range("A200:U200")
$headerRange.Copy()
$styleBackup.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$headerRange.ClearFormats()

$dataRange = $ws.Range("A1:U79")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

$styleBackup.Copy()
$headerRange.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$styleBackup.ClearContents()
$styleBackup.ClearFormats()
$excel.CutCopyMode = $false

# --- 3) Freeze the header row ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Host "Edit applied"
